$d = $word.ActiveDocument

# --- 1) Recolour the two "ready to start" (orange) items that are now blockers (red) ---
# "Rough design layout " run: FFC000 -> FF0000
$rng1 = $d.Content.Duplicate
$null = $rng1.Find.Execute("Rough design layout ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng1.Font.Color = 255

# "Usability testing " run (first occurrence, the one under "Rough design layout"): FFC000 -> FF0000
$rng2 = $d.Content.Duplicate
$null = $rng2.Find.Execute("Usability testing ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng2.Font.Color = 255

# --- 2) Move the "_GoBack" bookmark from the end of the "Enter production." paragraph ---
#        down to the blank paragraph at the very end of the document (after the new
#        "Created rough personas..." / Green / Orange / Red key paragraphs). ---
$old = $d.Bookmarks.Item("_GoBack")
$old.Delete()

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$markRange = $lastPara.Range.Duplicate
$markRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $markRange)
